$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three new rows before the current row 254, shifting existing
# rows 254:285 down to 257:288 (dimension grows from R285 to R288).
$ws.Rows.Item(254).Insert()
$ws.Rows.Item(255).Insert()
$ws.Rows.Item(256).Insert()

# New row 254: Ají, Americana (o), Primera
$ws.Cells.Item(254, 1).Value = 2
$ws.Cells.Item(254, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(254, 3).Value = "Coquimbo"
$ws.Cells.Item(254, 4).Value = 44748
$ws.Cells.Item(254, 5).Value = 4
$ws.Cells.Item(254, 6).Value = 100112021
$ws.Cells.Item(254, 7).Value = "Ají"
$ws.Cells.Item(254, 8).Value = "Americana (o)"
$ws.Cells.Item(254, 9).Value = "Primera"
$ws.Cells.Item(254, 10).Value = 200
$ws.Cells.Item(254, 11).Value = 30000
$ws.Cells.Item(254, 12).Value = 33000
$ws.Cells.Item(254, 13).Value = 31500
$ws.Cells.Item(254, 14).Value = "`$/caja 25 kilos"
$ws.Cells.Item(254, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(254, 16).Value = 1260
$ws.Cells.Item(254, 17).Value = 25
$ws.Cells.Item(254, 18).Value = "Hortaliza"

# New row 255: Ají, Americana (o), Segunda
$ws.Cells.Item(255, 1).Value = 2
$ws.Cells.Item(255, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(255, 3).Value = "Coquimbo"
$ws.Cells.Item(255, 4).Value = 44748
$ws.Cells.Item(255, 5).Value = 4
$ws.Cells.Item(255, 6).Value = 100112021
$ws.Cells.Item(255, 7).Value = "Ají"
$ws.Cells.Item(255, 8).Value = "Americana (o)"
$ws.Cells.Item(255, 9).Value = "Segunda"
$ws.Cells.Item(255, 10).Value = 100
$ws.Cells.Item(255, 11).Value = 20000
$ws.Cells.Item(255, 12).Value = 23000
$ws.Cells.Item(255, 13).Value = 21500
$ws.Cells.Item(255, 14).Value = "`$/caja 25 kilos"
$ws.Cells.Item(255, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(255, 16).Value = 860
$ws.Cells.Item(255, 17).Value = 25
$ws.Cells.Item(255, 18).Value = "Hortaliza"

# New row 256: Ají, Inferno, Primera
$ws.Cells.Item(256, 1).Value = 2
$ws.Cells.Item(256, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(256, 3).Value = "Coquimbo"
$ws.Cells.Item(256, 4).Value = 44748
$ws.Cells.Item(256, 5).Value = 4
$ws.Cells.Item(256, 6).Value = 100112021
$ws.Cells.Item(256, 7).Value = "Ají"
$ws.Cells.Item(256, 8).Value = "Inferno"
$ws.Cells.Item(256, 9).Value = "Primera"
$ws.Cells.Item(256, 10).Value = 160
$ws.Cells.Item(256, 11).Value = 23000
$ws.Cells.Item(256, 12).Value = 25000
$ws.Cells.Item(256, 13).Value = 24000
$ws.Cells.Item(256, 14).Value = "`$/caja 25 kilos"
$ws.Cells.Item(256, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(256, 16).Value = 960
$ws.Cells.Item(256, 17).Value = 25
$ws.Cells.Item(256, 18).Value = "Hortaliza"
